# EnrollMultipleLearners.xlsx — the startdate/enddate columns were re-entered
# as plain text (e.g. "2021-09-07") instead of real Excel dates, and the
# header row was given that same "Text" number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the header row (A1:D1) and the date columns (C2:D5) a Text number
# format ("@" -> numFmtId 49) before typing the new literal date strings,
# so Excel stores them as text rather than re-parsing them as serial dates.
$ws.Range("A1:D1").NumberFormat = "@"
$ws.Range("C2:D5").NumberFormat = "@"

$ws.Range("C2").Value = "2021-09-07"
$ws.Range("C3").Value = "2021-09-07"
$ws.Range("C4").Value = "2021-09-07"
$ws.Range("C5").Value = "2021-09-07"

$ws.Range("D2").Value = "2022-06-07"
$ws.Range("D3").Value = "2022-06-07"
$ws.Range("D4").Value = "2022-06-07"
$ws.Range("D5").Value = "2022-06-07"

# Move/record the active selection.
$ws.Range("C9").Select()

# Touch the sheet's page setup (orientation) so the print/page-setup block
# is (re)written for the sheet, as happened in the original edit.
$ws.PageSetup.Orientation = 1
